# Apply each division-problem replacement to its specific table cell.
# Using Replace=1 (wdReplaceOne) scoped to the cell Range avoids touching
# other cells that may share identical text (e.g. duplicate "75÷8=").
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$cell = $t.Cell(1, 1).Range
$cell.Find.Execute("43÷8=", $true, $false, $false, $false, $false, $true, 0, $false, "21÷5=", 1) | Out-Null
$cell = $t.Cell(1, 2).Range
$cell.Find.Execute("91÷4=", $true, $false, $false, $false, $false, $true, 0, $false, "17÷2=", 1) | Out-Null
$cell = $t.Cell(1, 3).Range
$cell.Find.Execute("51÷2=", $true, $false, $false, $false, $false, $true, 0, $false, "36÷5=", 1) | Out-Null
$cell = $t.Cell(1, 4).Range
$cell.Find.Execute("72÷9=", $true, $false, $false, $false, $false, $true, 0, $false, "34÷4=", 1) | Out-Null
$cell = $t.Cell(1, 5).Range
$cell.Find.Execute("93÷7=", $true, $false, $false, $false, $false, $true, 0, $false, "66÷5=", 1) | Out-Null
$cell = $t.Cell(5, 1).Range
$cell.Find.Execute("16÷7=", $true, $false, $false, $false, $false, $true, 0, $false, "99÷3=", 1) | Out-Null
$cell = $t.Cell(5, 2).Range
$cell.Find.Execute("12÷5=", $true, $false, $false, $false, $false, $true, 0, $false, "81÷5=", 1) | Out-Null
$cell = $t.Cell(5, 3).Range
$cell.Find.Execute("84÷5=", $true, $false, $false, $false, $false, $true, 0, $false, "62÷4=", 1) | Out-Null
$cell = $t.Cell(5, 4).Range
$cell.Find.Execute("75÷8=", $true, $false, $false, $false, $false, $true, 0, $false, "25÷3=", 1) | Out-Null
$cell = $t.Cell(5, 5).Range
$cell.Find.Execute("75÷8=", $true, $false, $false, $false, $false, $true, 0, $false, "20÷7=", 1) | Out-Null
$cell = $t.Cell(9, 1).Range
$cell.Find.Execute("29÷9=", $true, $false, $false, $false, $false, $true, 0, $false, "56÷6=", 1) | Out-Null
$cell = $t.Cell(9, 2).Range
$cell.Find.Execute("85÷4=", $true, $false, $false, $false, $false, $true, 0, $false, "21÷4=", 1) | Out-Null
$cell = $t.Cell(9, 3).Range
$cell.Find.Execute("95÷6=", $true, $false, $false, $false, $false, $true, 0, $false, "36÷6=", 1) | Out-Null
$cell = $t.Cell(9, 4).Range
$cell.Find.Execute("52÷8=", $true, $false, $false, $false, $false, $true, 0, $false, "37÷7=", 1) | Out-Null
$cell = $t.Cell(9, 5).Range
$cell.Find.Execute("29÷8=", $true, $false, $false, $false, $false, $true, 0, $false, "66÷7=", 1) | Out-Null
$cell = $t.Cell(13, 1).Range
$cell.Find.Execute("89÷5=", $true, $false, $false, $false, $false, $true, 0, $false, "29÷5=", 1) | Out-Null
$cell = $t.Cell(13, 2).Range
$cell.Find.Execute("76÷3=", $true, $false, $false, $false, $false, $true, 0, $false, "99÷2=", 1) | Out-Null
$cell = $t.Cell(13, 3).Range
$cell.Find.Execute("96÷9=", $true, $false, $false, $false, $false, $true, 0, $false, "17÷5=", 1) | Out-Null
$cell = $t.Cell(13, 4).Range
$cell.Find.Execute("31÷9=", $true, $false, $false, $false, $false, $true, 0, $false, "26÷3=", 1) | Out-Null
$cell = $t.Cell(13, 5).Range
$cell.Find.Execute("53÷5=", $true, $false, $false, $false, $false, $true, 0, $false, "60÷7=", 1) | Out-Null
$cell = $t.Cell(17, 1).Range
$cell.Find.Execute("87÷6=", $true, $false, $false, $false, $false, $true, 0, $false, "73÷4=", 1) | Out-Null
$cell = $t.Cell(17, 2).Range
$cell.Find.Execute("74÷4=", $true, $false, $false, $false, $false, $true, 0, $false, "53÷2=", 1) | Out-Null
$cell = $t.Cell(17, 3).Range
$cell.Find.Execute("60÷2=", $true, $false, $false, $false, $false, $true, 0, $false, "12÷5=", 1) | Out-Null
$cell = $t.Cell(17, 4).Range
$cell.Find.Execute("64÷8=", $true, $false, $false, $false, $false, $true, 0, $false, "61÷5=", 1) | Out-Null
$cell = $t.Cell(17, 5).Range
$cell.Find.Execute("54÷3=", $true, $false, $false, $false, $false, $true, 0, $false, "39÷3=", 1) | Out-Null
